$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1815789473684211
$ws.Range("C2").Value = 0.5605263157894737
$ws.Range("J2").Value = 0.01052631578947368
$ws.Range("P2").Value = 0.1342105263157895
$ws.Range("S2").Value = 0.1131578947368421
$ws.Range("B3").Value = 0.009049773755656109
$ws.Range("C3").Value = 0.03167420814479638
$ws.Range("J3").Value = 0.03167420814479638
$ws.Range("P3").Value = 0.755656108597285
$ws.Range("S3").Value = 0.1719457013574661
$ws.Range("J4").Value = 0.08163265306122448
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2040816326530612
$ws.Range("B6").Value = 0.06134969325153374
$ws.Range("D6").Value = 0.006134969325153374
$ws.Range("F6").Value = 0.03987730061349693
$ws.Range("J6").Value = 0.2975460122699387
$ws.Range("O6").Value = 0.01533742331288344
$ws.Range("Q6").Value = 0.1871165644171779
$ws.Range("R6").Value = 0.06134969325153374
$ws.Range("S6").Value = 0.3312883435582822
$ws.Range("B7").Value = 0.09965635738831616
$ws.Range("D7").Value = 0.0274914089347079
$ws.Range("F7").Value = 0.0584192439862543
$ws.Range("J7").Value = 0.1443298969072165
$ws.Range("O7").Value = 0.0274914089347079
$ws.Range("Q7").Value = 0.1718213058419244
$ws.Range("R7").Value = 0.06529209621993128
$ws.Range("S7").Value = 0.4054982817869416
$ws.Range("B8").Value = 0.08811188811188811
$ws.Range("D8").Value = 0.01258741258741259
$ws.Range("F8").Value = 0.08111888111888112
$ws.Range("J8").Value = 0.0979020979020979
$ws.Range("O8").Value = 0.01818181818181818
$ws.Range("Q8").Value = 0.1972027972027972
$ws.Range("R8").Value = 0.07132867132867132
$ws.Range("S8").Value = 0.4335664335664335
$ws.Range("B9").Value = 0.08813559322033898
$ws.Range("D9").Value = 0.01016949152542373
$ws.Range("F9").Value = 0.04067796610169491
$ws.Range("J9").Value = 0.1254237288135593
$ws.Range("O9").Value = 0.02711864406779661
$ws.Range("Q9").Value = 0.1864406779661017
$ws.Range("R9").Value = 0.06440677966101695
$ws.Range("S9").Value = 0.4576271186440678
$ws.Range("B10").Value = 0.1027692307692308
$ws.Range("D10").Value = 0.01784615384615384
$ws.Range("E10").Value = 0.002461538461538462
$ws.Range("F10").Value = 0.07630769230769231
$ws.Range("J10").Value = 0.1089230769230769
$ws.Range("O10").Value = 0.01292307692307692
$ws.Range("Q10").Value = 0.2098461538461538
$ws.Range("R10").Value = 0.07261538461538461
$ws.Range("S10").Value = 0.3963076923076923
$ws.Range("G11").Value = 0.1318681318681319
$ws.Range("J11").Value = 0.06043956043956044
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.6510989010989011
$ws.Range("S11").Value = 0.01373626373626374
$ws.Range("G12").Value = 0.7660377358490567
$ws.Range("J12").Value = 0.1735849056603773
$ws.Range("K12").Value = 0.007547169811320755
$ws.Range("L12").Value = 0.02641509433962264
$ws.Range("S12").Value = 0.02641509433962264
$ws.Range("G13").Value = 0.6883116883116883
$ws.Range("J13").Value = 0.2987012987012987
$ws.Range("S13").Value = 0.01298701298701299
$ws.Range("F15").Value = 0.02236421725239617
$ws.Range("H15").Value = 0.1821086261980831
$ws.Range("I15").Value = 0.08626198083067092
$ws.Range("J15").Value = 0.3162939297124601
$ws.Range("K15").Value = 0.07987220447284345
$ws.Range("M15").Value = 0.02236421725239617
$ws.Range("O15").Value = 0.04792332268370607
$ws.Range("S15").Value = 0.2428115015974441
$ws.Range("F16").Value = 0.0199203187250996
$ws.Range("H16").Value = 0.2231075697211155
$ws.Range("I16").Value = 0.09561752988047809
$ws.Range("J16").Value = 0.3306772908366534
$ws.Range("K16").Value = 0.09561752988047809
$ws.Range("M16").Value = 0.05179282868525897
$ws.Range("O16").Value = 0.09163346613545817
$ws.Range("S16").Value = 0.09163346613545817
$ws.Range("F17").Value = 0.02160493827160494
$ws.Range("H17").Value = 0.2098765432098765
$ws.Range("I17").Value = 0.08796296296296297
$ws.Range("J17").Value = 0.3796296296296297
$ws.Range("K17").Value = 0.09722222222222222
$ws.Range("M17").Value = 0.01697530864197531
$ws.Range("O17").Value = 0.07098765432098765
$ws.Range("S17").Value = 0.1157407407407407
$ws.Range("F18").Value = 0.01339285714285714
$ws.Range("H18").Value = 0.2276785714285714
$ws.Range("I18").Value = 0.05357142857142857
$ws.Range("J18").Value = 0.3616071428571428
$ws.Range("K18").Value = 0.15625
$ws.Range("M18").Value = 0.03125
$ws.Range("N18").Value = 0.004464285714285714
$ws.Range("O18").Value = 0.07589285714285714
$ws.Range("S18").Value = 0.07589285714285714
$ws.Range("F19").Value = 0.02259887005649718
$ws.Range("H19").Value = 0.2310734463276836
$ws.Range("I19").Value = 0.09717514124293786
$ws.Range("J19").Value = 0.3581920903954802
$ws.Range("K19").Value = 0.08757062146892655
$ws.Range("M19").Value = 0.02259887005649718
$ws.Range("O19").Value = 0.07062146892655367
$ws.Range("S19").Value = 0.1101694915254237
